$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# New handback identifiers / timestamps (per commit: "Generate Report for
# Handback" - a fresh handback run produced new handoff ids and new
# handoff/handback datetimes).
# ---------------------------------------------------------------------------
$oldMdA = "2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$newMdA = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.md"

$oldMdB = "e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"
$newMdB = "ffff1b7f35c7-1b13-4b8b-9b9b-e638c7ed79da.md"

$newZhXlf = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.zh-cn.xlf"
$newDeXlf = "567f5b35-998d-4faf-aa16-a2c0f7cb7d55.d0dc26581f56af7f0997c05bfe685f1717bc0c37.de-de.xlf"

$newZhHandoffTime = "2016-03-21 15:03:59"
$newZhHandbackTime = "2016-03-21 15:04:25"
$newDeHandoffTime = "2016-03-21 15:04:04"
$newDeHandbackTime = "2016-03-21 15:04:31"

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$wsOverview = $wb.Worksheets.Item("Overview")

$urlOverviewA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$urlOverviewA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"

$wsOverview.Range("A2").Value = $newMdA
$wsOverview.Range("A3").Value = $newMdB

$wsOverview.UsedRange.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $urlOverviewA2, "", "", $newMdA)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $urlOverviewA3, "", "", $newMdB)

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$wsZh = $wb.Worksheets.Item("zh-cn")

$urlZhA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$urlZhD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c4d779ba04829f60b44c18ecce367a5e341734c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.zh-cn.xlf"
$urlZhF2 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0b31f8556abf8d1f233da753bcbfbe7583bfe77e/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$urlZhG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/844aa5308dfc5272a86ea64056e8e01d14eb8763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.zh-cn.xlf"
$urlZhA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"
$urlZhD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c4d779ba04829f60b44c18ecce367a5e341734c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.zh-cn.xlf"
$urlZhF3 = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/0b31f8556abf8d1f233da753bcbfbe7583bfe77e/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"
$urlZhG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/844aa5308dfc5272a86ea64056e8e01d14eb8763/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.zh-cn.xlf"

# Row 2
$wsZh.Range("A2").Value = $newMdA
$wsZh.Range("D2").Value = $newZhXlf
$wsZh.Range("E2").Value = $newZhHandoffTime
$wsZh.Range("F2").Value = $newMdA
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandbackTime

# Row 3
$wsZh.Range("A3").Value = $newMdB
$wsZh.Range("D3").Value = $newZhXlf
$wsZh.Range("E3").Value = $newZhHandoffTime
$wsZh.Range("F3").Value = $newMdB
$wsZh.Range("G3").Value = $newZhXlf
$wsZh.Range("H3").Value = $newZhHandbackTime

$wsZh.UsedRange.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $urlZhA2, "", "", $newMdA)
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $urlZhD2, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $urlZhF2, "", "", $newMdA)
$wsZh.Hyperlinks.Add($wsZh.Range("G2"), $urlZhG2, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $urlZhA3, "", "", $newMdB)
$wsZh.Hyperlinks.Add($wsZh.Range("D3"), $urlZhD3, "", "", $newZhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $urlZhF3, "", "", $newMdB)
$wsZh.Hyperlinks.Add($wsZh.Range("G3"), $urlZhG3, "", "", $newZhXlf)

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$wsDe = $wb.Worksheets.Item("de-de")

$urlDeA2 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$urlDeD2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f87811625197db4aa2d219e29072069f73769b2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.de-de.xlf"
$urlDeF2 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/01b5a4e0b5f0fdef90b4ecf9918e035a3c8d89ce/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md"
$urlDeG2 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d52a0eab2f9310c7209f9e1a2296c0f78b8fe03f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.de-de.xlf"
$urlDeA3 = "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"
$urlDeD3 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f87811625197db4aa2d219e29072069f73769b2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.de-de.xlf"
$urlDeF3 = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/01b5a4e0b5f0fdef90b4ecf9918e035a3c8d89ce/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md"
$urlDeG3 = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/d52a0eab2f9310c7209f9e1a2296c0f78b8fe03f/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.de-de.xlf"

# Row 2
$wsDe.Range("A2").Value = $newMdA
$wsDe.Range("D2").Value = $newDeXlf
$wsDe.Range("E2").Value = $newDeHandoffTime
$wsDe.Range("F2").Value = $newMdA
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newDeHandbackTime

# Row 3
$wsDe.Range("A3").Value = $newMdB
$wsDe.Range("D3").Value = $newDeXlf
$wsDe.Range("E3").Value = $newDeHandoffTime
$wsDe.Range("F3").Value = $newMdB
$wsDe.Range("G3").Value = $newDeXlf
$wsDe.Range("H3").Value = $newDeHandbackTime

$wsDe.UsedRange.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $urlDeA2, "", "", $newMdA)
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $urlDeD2, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $urlDeF2, "", "", $newMdA)
$wsDe.Hyperlinks.Add($wsDe.Range("G2"), $urlDeG2, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $urlDeA3, "", "", $newMdB)
$wsDe.Hyperlinks.Add($wsDe.Range("D3"), $urlDeD3, "", "", $newDeXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $urlDeF3, "", "", $newMdB)
$wsDe.Hyperlinks.Add($wsDe.Range("G3"), $urlDeG3, "", "", $newDeXlf)
